$wb = $excel.ActiveWorkbook

# --- Metadata sheet updates ---
$meta = $wb.Worksheets.Item("Metadata")

# Version: 2.1.0 -> 2.2.0-ballot
$meta.Range("B3").Value = "2.2.0-ballot"

# Date: 2025-12-19T08:22:07+00:00 -> 2025-12-19T09:47:21+00:00
$meta.Range("B8").Value = "2025-12-19T09:47:21+00:00"

# Base Definition: append the |4.0.1 version tag
$meta.Range("B18").Value = "http://hl7.org/fhir/StructureDefinition/Extension|4.0.1"

# --- Elements sheet updates ---
$elem = $wb.Worksheets.Item("Elements")

# Extension.value[x] Type(s) cell: append |2.2.0-ballot to the referenced profile, keep trailing newline
$elem.Range("K6").Value = "Reference(https://interop.esante.gouv.fr/ig/fhir/tddui/StructureDefinition/tddui-careplan-projet-personnalise|2.2.0-ballot)`n"

# Column K (11) grew wider to fit the longer text (bestFit recalculation): 86.23828125 -> ~95.746
$elem.Columns.Item(11).ColumnWidth = 94.8
